$wb = $excel.ActiveWorkbook

# Insert the new worksheet right before "Small N+zero+names" so it lands
# between "Small N w zero BV" and "Small N+zero+names".
$target = $wb.Worksheets.Item("Small N+zero+names")
$ws = $wb.Worksheets.Add($target)
$ws.Name = "No data"

# Staircase of words across the sheet.
$ws.Range("A1").Value = "There"
$ws.Range("B1").Value = "is"
$ws.Range("C1").Value = "no"
$ws.Range("A2").Value = "data"
$ws.Range("B3").Value = "in"
$ws.Range("C4").Value = "here"
$ws.Range("D5").Value = "hmm"

# Match the selection / active-sheet state captured in the workbook.
$ws.Range("D6").Select()
$ws.Activate()
